$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "student_id" column (I) for each student row (2-22). Base the
# formatting on the existing roster data cells (Book Antiqua 11, black)
# and center it - this mirrors the new cell style introduced alongside
# the data.
$ws.Range("B2").Copy()
$idRange = $ws.Range("I2:I22")
$idRange.PasteSpecial(-4122)
$idRange.HorizontalAlignment = -4108
$excel.CutCopyMode = $false

$ws.Range("I2").Value = "18-0198"
$ws.Range("I3").Value = "18-0219"
$ws.Range("I4").Value = "18-0234"
$ws.Range("I5").Value = "18-0202"
$ws.Range("I6").Value = "18-0018"
$ws.Range("I7").Value = "18-0141"
$ws.Range("I8").Value = "18-0139"
$ws.Range("I9").Value = "18-0096"
$ws.Range("I10").Value = "18-0079"
$ws.Range("I11").Value = "18-0091"
$ws.Range("I12").Value = "18-0067"
$ws.Range("I13").Value = "18-0076"
$ws.Range("I14").Value = "18-0077"
$ws.Range("I15").Value = "18-0015"
$ws.Range("I16").Value = "18-0231"
$ws.Range("I17").Value = "18-0080"
$ws.Range("I18").Value = "18-0140"
$ws.Range("I19").Value = "18-0150"
$ws.Range("I20").Value = "18-0183"
$ws.Range("I21").Value = "18-0090"
$ws.Range("I22").Value = "18-0237"

# Match the author's final on-screen selection/scroll position.
$ws.Range("A19:E20").Select()
